$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 81

# Columns A-D hold text-like values (date/time/weekday/week as strings in the
# source data). Force text interpretation so Excel doesn't auto-convert the
# date/time strings into serial numbers or drop the leading zero in "07".
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-02-21"
$ws.Cells.Item($row, 2).Value = "09:45:27"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).Value = "07"

# Restore the default cell style so the new row matches the unstyled look of
# the other data rows (only header row uses an explicit style).
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 130043
$ws.Cells.Item($row, 6).Value = 140952
$ws.Cells.Item($row, 7).Value = 171768
$ws.Cells.Item($row, 8).Value = 154225
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146181
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193110
$ws.Cells.Item($row, 14).Value = 115299
$ws.Cells.Item($row, 15).Value = 46011
$ws.Cells.Item($row, 16).Value = 29179
$ws.Cells.Item($row, 17).Value = 67662
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47507
$ws.Cells.Item($row, 20).Value = -1
